$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.405.21'
$ws.Range("E2").Value = '  -1.61%  '
$ws.Range("D3").Value = '3.397.68'
$ws.Range("E3").Value = '  -0.45%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '566.68'
$ws.Range("E5").Value = '  -0.95%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '156.08'
$ws.Range("E6").Value = '  -0.36%  '
$ws.Range("D8").Value = '3.398.16'
$ws.Range("E8").Value = '  -0.55%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.565'
$ws.Range("E9").Value = '  -8.31%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.25'
$ws.Range("E10").Value = '  +1.04%  '
$ws.Range("E11").Value = '  -3.30%  '
$ws.Range("E12").Value = '  -4.11%  '
$ws.Range("E13").Value = '  -0.56%  '
$ws.Range("E14").Value = '  -0.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.84'
$ws.Range("E15").Value = '  -4.12%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000170'
$ws.Range("E16").Value = '  -9.13%  '
$ws.Range("D17").Value = '63.484.58'
$ws.Range("E17").Value = '  -1.55%  '
$ws.Range("D18").Value = '3.403.34'
$ws.Range("E18").Value = '  -0.60%  '
$ws.Range("E19").Value = '  -4.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.48'
$ws.Range("E20").Value = '  -3.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '383.36'
$ws.Range("E21").Value = '  +2.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.72'
$ws.Range("E22").Value = '  -3.39%  '
$ws.Range("E23").Value = '  +0.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.04'
$ws.Range("E24").Value = '  -1.96%  '
$ws.Range("E25").Value = '  -7.21%  '
$ws.Range("E26").Value = '  -3.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.67'
$ws.Range("E27").Value = '  -5.69%  '
$ws.Range("E28").Value = '  +0.41%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("E30").Value = '  -2.40%  '
$ws.Range("E31").Value = '  -7.38%  '
$ws.Range("E32").Value = '  -2.54%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '22.80'
$ws.Range("E33").Value = '  -1.22%  '
$ws.Range("E34").Value = '  -4.34%  '
$ws.Range("E35").Value = '  -6.96%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '161.03'
$ws.Range("E36").Value = '  +0.37%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.838'
$ws.Range("E37").Value = '  +9.03%  '
$ws.Range("E38").Value = '  -4.15%  '
$ws.Range("D39").Value = '2.807.98'
$ws.Range("E39").Value = '  -1.30%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '25.85'
$ws.Range("E40").Value = '  -2.87%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '42.87'
$ws.Range("E41").Value = '  -0.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0716'
$ws.Range("E42").Value = '  -5.79%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.37'
$ws.Range("E43").Value = '  -7.44%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '25.57'
$ws.Range("E44").Value = '  -3.50%  '
$ws.Range("E45").Value = '  -5.76%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0302'
$ws.Range("E46").Value = '  -3.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '325.82'
$ws.Range("E47").Value = '  +1.99%  '
$ws.Range("E48").Value = '  +8.08%  '
$ws.Range("E49").Value = '  -5.28%  '
$ws.Range("E50").Value = '  -5.55%  '
$ws.Range("E51").Value = '  -4.90%  '
